$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old 2000年-2009年 rows (old rows 2..11) so that the existing
# 2010年-2015年 rows (old rows 12..17) shift up to become rows 2..7.
$ws.Range("A2:D11").EntireRow.Delete()
